$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = "https://selectorshub.com/xpath-practice-page/"
$ws.Range("H2").Value = 6000

# --- Row 4 (write new strings first so shared-string table order matches target) ---
$ws.Range("D4").Value = "Username"
$ws.Range("E4").Value = "a"

# --- Row 3 ---
$ws.Range("D3").Value = "Users Table"
$ws.Range("E3").Value = "h3"
$ws.Range("C3").Value = "scroll"

# --- Row 5 ---
$ws.Range("C5").Value = "select"
$ws.Range("D5").Value = "checkbox"
$ws.Range("E5").Value = "input"

# --- Row 4 remaining ---
$ws.Range("C4").Value = "waitfortext"

# Remove the mailto hyperlink on F4 and clear its contents entirely
$ws.Range("F4").Hyperlinks.Delete()
$ws.Range("F4").Clear()

# F3 becomes a normal-styled numeric cell (was an empty Hyperlink-styled cell)
$ws.Range("G3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = 1

# F5 and H5 are removed entirely
$ws.Range("F5").Clear()
$ws.Range("H5").Clear()

# Rows 6-11 are removed entirely without shifting rows below them
$ws.Rows("6:11").Clear()

# Selection moves to the whole of row 5
$ws.Rows("5:5").Select() | Out-Null

"done"
